# The deck ships with a custom "Integral" theme (ppt/theme/theme1.xml,
# linked from the slide master) and, separately, a stock "Office Theme"
# that is only used by the notes master (ppt/theme/theme2.xml).
#
# This edit switches the presentation's applied design back to the
# default Office theme colors by rewriting every slot of the active
# 12-colour theme scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) -
# the font scheme and format (fill/line/effect) scheme are already
# identical between the two themes, so the colours are the only thing
# that needs to change to go from "Integral" to "Office Theme".

$p = $ppt.ActivePresentation

# Office theme colour scheme, converted to VBA RGB() integers
# (R + G*256 + B*65536) in clrScheme order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

# The theme colour scheme is shared across the whole deck, so editing it
# through any slide updates the slide master's theme part for everyone.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeThemeColors[$i - 1]
}
